$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order lines (SKU, Name, Quantity, Cost Per, Total Cost) appended below
# the existing table, which ran through row 18.
$data = @(
    ,@('7520950', 'Can - Pizza Sauce', '2', '41.22', '82.44')
    ,@('8457368', 'Oil - Corn', '4', '35.30', '141.20')
    ,@('8255796', 'Tuna White Chunk (Pouch)', '5', '72.00', '360.00')
    ,@('1028165', 'Pickle - Dill Chip', '2', '34.09', '68.18')
    ,@('7722184', 'Parmesan (Grated)', '1', '59.95', '59.95')
    ,@('1027629', 'Cheddar - (Sliced)', '12', '34.91', '418.92')
    ,@('1035842', 'Feta - Pail', '1', '92.87', '92.87')
    ,@('6364494', 'Yogurt - Greek (Bulk)', '2', '26.96', '53.92')
    ,@('1365278', 'Vegan Chicken Tenders', '2', '87.80', '175.60')
    ,@('3275539', 'Sauerkraut', '2', '19.35', '38.70')
    ,@('1028188', 'Tortellini - Cheese', '4', '32.22', '128.88')
    ,@('7529232', 'Wrap - Wheat (10")', '1', '31.32', '31.32')
    ,@('2825368', 'Sausage - Chicken Patty', '4', '50.01', '200.04')
    ,@('4157160', 'Spanakopita', '2', '75.51', '151.02')
    ,@('2477933', 'Bacon (Pre-Cooked)', '20', '36.40', '728.00')
    ,@('9546982', 'Arugula - Fresh', '4', '20.35', '81.40')
    ,@('2054542', 'Carrots - Jumbo Fresh', '1', '44.13', '44.13')
)

$startRow = 19
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 5))

    # The existing rows store every column -- including the numeric-looking
    # ones -- as plain text. Flip to a text number format before assigning so
    # Excel does not coerce these numeric-looking strings into real numbers,
    # then restore the default style so no extra formatting is left behind.
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]

    $rowRange.Style = "Normal"
}

Write-Host "Done adding rows 19-35"
